$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.5418
$ws.Range("C4").Value = -12.44220000000001
$ws.Range("B7").Value = 5.632200000000001
$ws.Range("A8").Value = -22.39810000000002
$ws.Range("A10").Value = -21.7641
$ws.Range("D10").Value = -7.595200000000001
$ws.Range("C11").Value = -12.42679999999999
$ws.Range("A12").Value = -21.57230000000001
$ws.Range("D12").Value = -7.303999999999993
$ws.Range("D13").Value = -8.974899999999993
$ws.Range("B14").Value = 6.068900000000006
$ws.Range("C14").Value = -12.3083
$ws.Range("D14").Value = -7.204900000000001
$ws.Range("B15").Value = 4.664599999999994
$ws.Range("A18").Value = -21.8847
$ws.Range("B18").Value = 6.358000000000004
$ws.Range("C18").Value = -12.0359
$ws.Range("C19").Value = -11.5629
$ws.Range("B20").Value = 8.843799999999998
$ws.Range("C21").Value = -12.3608
$ws.Range("A25").Value = -21.51099999999999
$ws.Range("C27").Value = -13.1947
$ws.Range("B29").Value = 4.942100000000006
$ws.Range("D29").Value = -7.350699999999996
$ws.Range("B30").Value = 5.324200000000001
$ws.Range("B31").Value = 4.391999999999999
$ws.Range("C31").Value = -13.60869999999999
$ws.Range("D32").Value = -8.960899999999993
$ws.Range("B35").Value = 8.205700000000002
$ws.Range("D35").Value = -8.225299999999994
$ws.Range("A37").Value = -19.14659999999999
$ws.Range("C38").Value = -12.9936
$ws.Range("B40").Value = 9.253599999999997
$ws.Range("C42").Value = -11.73010000000001
$ws.Range("D43").Value = -7.963899999999999
$ws.Range("B44").Value = 5.482800000000001
$ws.Range("C44").Value = -13.42779999999999
$ws.Range("C47").Value = -12.58
$ws.Range("D48").Value = -7.444999999999997
$ws.Range("D49").Value = -8.231000000000003
$ws.Range("B50").Value = 4.728699999999996
$ws.Range("D50").Value = -8.296499999999995
$ws.Range("D51").Value = -8.088899999999994
$ws.Range("B54").Value = 4.791599999999999
$ws.Range("A55").Value = -21.9429
$ws.Range("C56").Value = -13.68249999999999
$ws.Range("D56").Value = -8.443300000000001
$ws.Range("C58").Value = -12.456
$ws.Range("D61").Value = -7.987699999999998
$ws.Range("C65").Value = -11.99810000000001
$ws.Range("A68").Value = -21.49410000000001
$ws.Range("B68").Value = 4.441600000000001
$ws.Range("D69").Value = -7.291499999999996
$ws.Range("D71").Value = -7.717399999999994
$ws.Range("C73").Value = -12.6933
$ws.Range("B76").Value = 6.043499999999994
$ws.Range("A77").Value = -20.13289999999998
$ws.Range("A78").Value = -20.06649999999998
$ws.Range("A79").Value = -19.96609999999998
$ws.Range("D79").Value = -6.361300000000004
$ws.Range("A80").Value = -19.81659999999997
$ws.Range("A81").Value = -21.7295
$ws.Range("D81").Value = -7.356299999999996
$ws.Range("A82").Value = -21.7775
$ws.Range("A84").Value = -22.046
$ws.Range("B87").Value = 4.790099999999995
$ws.Range("B88").Value = 4.703499999999996
$ws.Range("C90").Value = -12.8744
$ws.Range("B92").Value = 5.908399999999999
$ws.Range("C92").Value = -12.3508
$ws.Range("D92").Value = -6.594799999999999
$ws.Range("C94").Value = -10.12060000000001
$ws.Range("C95").Value = -12.34279999999999
$ws.Range("B96").Value = 5.226100000000003
$ws.Range("B98").Value = 6.225499999999998
$ws.Range("A101").Value = -21.57319999999998
$ws.Range("B101").Value = 6.0762
$ws.Range("C101").Value = -12.69669999999999
$ws.Range("A102").Value = -21.8016
$ws.Range("B102").Value = 5.624500000000003
